# auto commit Tue Dec 19 16:42:22 CET 2017
# Re-layout the benefits/drawbacks table: drop a redundant bullet, add a
# third "grouping" column (SPs / Users / Misc) that buckets each row, and
# shade the three groups. Also nudges the window position and selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the row/column layout is being substantially reshuffled.
$ws.Cells.Clear()

# ---------------------------------------------------------------------
# Cell values
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Benefits"
$ws.Range("B1").Value = "Drawbacks"

$ws.Range("A2").Value = "Low costs"
$ws.Range("B2").Value = "Large number of attack vectors"
$ws.Range("C2").Value = "SPs"

$ws.Range("A3").Value = "Easy to implement"
$ws.Range("B3").Value = "Anomaly detection costly"

$ws.Range("A4").Value = "Replaceable when compromised"
$ws.Range("B4").Value = "Attacks are simple to carry out"

$ws.Range("A5").Value = "Revokable by administrator"
$ws.Range("B5").Value = "Attack automation simple"

$ws.Range("A6").Value = "Enforceable policies"

$ws.Range("A7").Value = "Fast entry on desktops"
$ws.Range("B7").Value = "Memory overload from too many passwords"
$ws.Range("C7").Value = "Users"

$ws.Range("A8").Value = "Most users already familiarized"
$ws.Range("B8").Value = "Suboptimal coping strategies"

$ws.Range("A9").Value = "Easy to learn"
$ws.Range("B9").Value = "Weak passwords are a risk for users and SPs "

$ws.Range("A10").Value = "Sharable with others"
$ws.Range("B10").Value = "Stronger passwords difficult to memorize"

$ws.Range("A11").Value = "High degree of control / freedom"
$ws.Range("B11").Value = "Entry on mobile devices difficult"

$ws.Range("B12").Value = "Mastery difficult"

$ws.Range("B13").Value = "Disliked by many users / perceived as burden"

$ws.Range("A14").Value = "Idenpendent of identification"
$ws.Range("C14").Value = "Misc"

$ws.Range("A15").Value = "Adjustable security level"

# ---------------------------------------------------------------------
# Header row formatting (bold + rule under it)
# ---------------------------------------------------------------------
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$ws.Rows.Item(1).RowHeight = 17
$header.Borders.Item(8).LineStyle = 1
$header.Borders.Item(8).Weight = 2
$header.Borders.Item(9).LineStyle = 1
$header.Borders.Item(9).Weight = -4138

# ---------------------------------------------------------------------
# Group shading - SPs block (rows 2-6), light grey (theme 0, ~5% darker)
# ---------------------------------------------------------------------
$spsBody = $ws.Range("A2:B6")
$spsBody.Interior.ThemeColor = 2
$spsBody.Interior.TintAndShade = -0.05
$spsTag = $ws.Range("C2:C6")
$spsTag.Interior.ThemeColor = 2
$spsTag.Interior.TintAndShade = -0.05
$spsTag.Borders.Item(7).LineStyle = 1
$spsTag.HorizontalAlignment = -4108
$spsTag.VerticalAlignment = -4108
$spsTag.Orientation = 180
$spsTag.Merge()

# ---------------------------------------------------------------------
# Group shading - Users block (rows 7-13), light green (theme 9, 80% lighter)
# ---------------------------------------------------------------------
$usersBody = $ws.Range("A7:B13")
$usersBody.Interior.ThemeColor = 10
$usersBody.Interior.TintAndShade = 0.8
$usersTag = $ws.Range("C7:C13")
$usersTag.Interior.ThemeColor = 10
$usersTag.Interior.TintAndShade = 0.8
$usersTag.Borders.Item(7).LineStyle = 1
$usersTag.HorizontalAlignment = -4108
$usersTag.VerticalAlignment = -4108
$usersTag.Orientation = 180
$usersTag.Merge()

# ---------------------------------------------------------------------
# Group shading - Misc block (rows 14-15), light blue (theme 4, 80% lighter)
# ---------------------------------------------------------------------
$miscBody = $ws.Range("A14:B15")
$miscBody.Interior.ThemeColor = 5
$miscBody.Interior.TintAndShade = 0.8
$miscTag = $ws.Range("C14:C15")
$miscTag.Interior.ThemeColor = 5
$miscTag.Interior.TintAndShade = 0.8
$miscTag.Borders.Item(7).LineStyle = 1
$miscTag.HorizontalAlignment = -4108
$miscTag.VerticalAlignment = -4108
$miscTag.Orientation = 180
$miscTag.Merge()

# ---------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 460
